# "Add files via upload" — append prior-season rows (2020 down to 2014) to
# the first sheet's year column (A), then move the active selection to
# A10, just below the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$years = @(2020, 2019, 2018, 2017, 2016, 2015, 2014)
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $years[$i]
}

$ws.Range("A10").Select() | Out-Null
